$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18/19 swap: WrappedEther <-> Polkadot (rows exchange coin identity + updated price/volume)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "6.72"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "2.957.44"
$ws.Range("E19").Value = "  +0.86%  "

# Price / Volume(1h) updates for remaining rows
# NumberFormat "@" (text) is (re)applied on column D only where the Price value
# itself changes, so it keeps being stored as text (matching the source data,
# which uses "." as a thousands separator, e.g. "62.800.59") instead of being
# auto-coerced into a number by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.800.59"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.966.44"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.70"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.71"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.966.46"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("E10").Value = "  +5.03%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000235"
$ws.Range("E13").Value = "  +4.31%  "
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.453.83"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.678.26"
$ws.Range("E17").Value = "  +2.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "442.17"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.44"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.95"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.29"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("E26").Value = "  +0.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.14"
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.14"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("E31").Value = "  -5.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.73"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0888"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.65"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("E41").Value = "  -0.22%  "
$ws.Range("E42").Value = "  -4.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.284"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.06"
$ws.Range("E44").Value = "  -6.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.731.93"
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "136.13"
$ws.Range("E46").Value = "  +1.74%  "
$ws.Range("E47").Value = "  -1.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "364.14"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.98"
$ws.Range("E51").Value = "  -3.96%  "
